$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "dSF" (column F) values for the specified rows, per the
# repulled data / mean recalculation described in the commit message.
$updates = @{
    4  = 3
    5  = 1
    16 = 0
    20 = -1
    27 = -2
    31 = -1
    33 = 2
    34 = -4
    35 = -1
    38 = 2
    49 = -3
    51 = 1
    52 = 2
    53 = -3
    54 = -2
    59 = -1
    69 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
